$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2088.353
$ws.Range("I4").Value = 492.85715
$ws.Range("J4").Value = 3205.2
$ws.Range("K4").Value = 492.85715
$ws.Range("L4").Value = 3205.2
$ws.Range("M4").Value = -378.85715
$ws.Range("N4").Value = -3433.2
$ws.Range("H11").Value = 1745.25
$ws.Range("I11").Value = 1745.25
$ws.Range("K11").Value = 1745.25
$ws.Range("M11").Value = -1605.25
$ws.Range("H18").Value = 321.73334
$ws.Range("I18").Value = 321.73334
$ws.Range("K18").Value = 321.73334
$ws.Range("M18").Value = -37.73334
$ws.Range("H43").Value = 986.6667
$ws.Range("I43").Value = 420
$ws.Range("J43").Value = 1175.5555
$ws.Range("K43").Value = 420
$ws.Range("L43").Value = 1175.5555
$ws.Range("M43").Value = -351
$ws.Range("N43").Value = -1313.5555
$ws.Range("H74").Value = 4558.1665
$ws.Range("I74").Value = 4245.273
$ws.Range("K74").Value = 4245.273
$ws.Range("M74").Value = -3309.273
$ws.Range("H77").Value = 4558.1665
$ws.Range("I77").Value = 4245.273
$ws.Range("K77").Value = 21226.365
$ws.Range("M77").Value = -16546.365
$ws.Range("H100").Value = 1644.3684
$ws.Range("I100").Value = 1726.5834
$ws.Range("J100").Value = 1503.4286
$ws.Range("K100").Value = 1726.5834
$ws.Range("L100").Value = 1503.4286
$ws.Range("M100").Value = -1185.5834
$ws.Range("N100").Value = -2585.4286
$ws.Range("H112").Value = 5755.125
$ws.Range("I112").Value = 17745
$ws.Range("J112").Value = 1758.5
$ws.Range("K112").Value = 53235
$ws.Range("L112").Value = 5275.5
$ws.Range("M112").Value = -52127
$ws.Range("N112").Value = -7491.5
$ws.Range("H126").Value = 74800
$ws.Range("J126").Value = 74800
$ws.Range("L126").Value = 74800
$ws.Range("N126").Value = -84680
$ws.Range("H129").Value = 914.9729599999999
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 930.6667
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 2792.0001
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -12792.0001
$ws.Range("H132").Value = 927.7778
$ws.Range("I132").Value = 936.5217
$ws.Range("J132").Value = 877.5
$ws.Range("K132").Value = 2809.5651
$ws.Range("L132").Value = 2632.5
$ws.Range("M132").Value = -279.5650999999998
$ws.Range("N132").Value = -7692.5
$ws.Range("H137").Value = 683013.5600000001
$ws.Range("I137").Value = 4084.3572
$ws.Range("J137").Value = 954585.25
$ws.Range("K137").Value = 12253.0716
$ws.Range("L137").Value = 2863755.75
$ws.Range("M137").Value = -9703.071599999999
$ws.Range("N137").Value = -2868855.75
$ws.Range("H138").Value = 3168.264
$ws.Range("I138").Value = 1391
$ws.Range("J138").Value = 4002.4897
$ws.Range("K138").Value = 4173
$ws.Range("L138").Value = 12007.4691
$ws.Range("M138").Value = 967
$ws.Range("N138").Value = -22287.4691
$ws.Range("H141").Value = 3258.0264
$ws.Range("I141").Value = 2023.409
$ws.Range("J141").Value = 4955.625
$ws.Range("K141").Value = 6070.227000000001
$ws.Range("L141").Value = 14866.875
$ws.Range("M141").Value = -890.2270000000008
$ws.Range("N141").Value = -25226.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22186.275
$ws.Range("I32").Value = 24082.021
$ws.Range("K32").Value = 24082.021
$ws.Range("M32").Value = -23795.021
$ws.Range("H61").Value = 8097.5586
$ws.Range("I61").Value = 6115.6113
$ws.Range("J61").Value = 10327.25
$ws.Range("K61").Value = 6115.6113
$ws.Range("L61").Value = 10327.25
$ws.Range("M61").Value = -5903.6113
$ws.Range("N61").Value = -10751.25
$ws.Range("H62").Value = 40249
$ws.Range("J62").Value = 40249
$ws.Range("L62").Value = 40249
$ws.Range("N62").Value = -41497
$ws.Range("H65").Value = 40249
$ws.Range("J65").Value = 40249
$ws.Range("L65").Value = 120747
$ws.Range("N65").Value = -126987
$ws.Range("H74").Value = 3820.275
$ws.Range("I74").Value = 1362.3
$ws.Range("K74").Value = 1362.3
$ws.Range("M74").Value = -488.3
$ws.Range("H77").Value = 3820.275
$ws.Range("I77").Value = 1362.3
$ws.Range("K77").Value = 6811.5
$ws.Range("M77").Value = -2443.5
$ws.Range("H97").Value = 771.43243
$ws.Range("I97").Value = 676.9655
$ws.Range("J97").Value = 1113.875
$ws.Range("K97").Value = 676.9655
$ws.Range("L97").Value = 1113.875
$ws.Range("M97").Value = -180.9655
$ws.Range("N97").Value = -2105.875
$ws.Range("H110").Value = 1972.5
$ws.Range("I110").Value = 1972.5
$ws.Range("K110").Value = 1972.5
$ws.Range("M110").Value = 72.5
$ws.Range("H136").Value = 8097.5586
$ws.Range("I136").Value = 6115.6113
$ws.Range("J136").Value = 10327.25
$ws.Range("K136").Value = 18346.8339
$ws.Range("L136").Value = 30981.75
$ws.Range("M136").Value = -15796.8339
$ws.Range("N136").Value = -36081.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2396.15
$ws.Range("I107").Value = 2411.1
$ws.Range("J107").Value = 2381.2
$ws.Range("K107").Value = 2411.1
$ws.Range("L107").Value = 2381.2
$ws.Range("M107").Value = -491.0999999999999
$ws.Range("N107").Value = -6221.2
$ws.Range("H134").Value = 45267.78
$ws.Range("I134").Value = 1956.55
$ws.Range("J134").Value = 334009.34
$ws.Range("K134").Value = 5869.65
$ws.Range("L134").Value = 1002028.02
$ws.Range("M134").Value = -3334.65
$ws.Range("N134").Value = -1007098.02

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 24791.105
$ws.Range("I59").Value = 10104
$ws.Range("K59").Value = 10104
$ws.Range("M59").Value = -8959
$ws.Range("H107").Value = 1057.9546
$ws.Range("I107").Value = 1105.4166
$ws.Range("K107").Value = 1105.4166
$ws.Range("M107").Value = 814.5834
$ws.Range("H122").Value = 7178.591
$ws.Range("I122").Value = 2189.7334
$ws.Range("J122").Value = 17869
$ws.Range("K122").Value = 6569.2002
$ws.Range("L122").Value = 53607
$ws.Range("M122").Value = -4119.2002
$ws.Range("N122").Value = -58507
$ws.Range("H132").Value = 2656.9302
$ws.Range("I132").Value = 2270.4443
$ws.Range("J132").Value = 4644.5713
$ws.Range("K132").Value = 6811.3329
$ws.Range("L132").Value = 13933.7139
$ws.Range("M132").Value = -4281.3329
$ws.Range("N132").Value = -18993.7139

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3300
$ws.Range("J117").Value = 3300
$ws.Range("L117").Value = 9900
$ws.Range("N117").Value = -16784
$ws.Range("H131").Value = 1135.3281
$ws.Range("I131").Value = 1686.2727
$ws.Range("J131").Value = 1020.98114
$ws.Range("K131").Value = 5058.8181
$ws.Range("L131").Value = 3062.94342
$ws.Range("M131").Value = -18.81810000000041
$ws.Range("N131").Value = -13142.94342

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3851.125
$ws.Range("I7").Value = 4361.8
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 4361.8
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -4249.8
$ws.Range("N7").Value = -3224
$ws.Range("H40").Value = 3245.5
$ws.Range("I40").Value = 3045.682
$ws.Range("J40").Value = 3685.1
$ws.Range("K40").Value = 3045.682
$ws.Range("L40").Value = 3685.1
$ws.Range("M40").Value = -2909.682
$ws.Range("N40").Value = -3957.1
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null
$ws.Range("H93").Value = 1172.5
$ws.Range("I93").Value = 1296.6666
$ws.Range("K93").Value = 1296.6666
$ws.Range("M93").Value = -48.66660000000002
$ws.Range("H122").Value = 5875.35
$ws.Range("I122").Value = 5662.0293
$ws.Range("J122").Value = 7084.1665
$ws.Range("K122").Value = 16986.0879
$ws.Range("L122").Value = 21252.4995
$ws.Range("M122").Value = -14536.0879
$ws.Range("N122").Value = -26152.4995
$ws.Range("H123").Value = 52576.47
$ws.Range("J123").Value = 52576.47
$ws.Range("L123").Value = 52576.47
$ws.Range("N123").Value = -62376.47
$ws.Range("H126").Value = 3851.125
$ws.Range("I126").Value = 4361.8
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 13085.4
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -10615.4
$ws.Range("N126").Value = -13940
$ws.Range("H141").Value = 69950
$ws.Range("J141").Value = 69950
$ws.Range("L141").Value = 69950
$ws.Range("N141").Value = -80310

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 44111.09
$ws.Range("J123").Value = 44111.09
$ws.Range("L123").Value = 44111.09
$ws.Range("N123").Value = -53911.09
